$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 518; existing rows 518:543 shift down to 520:545.
$ws.Rows.Item(518).Insert()
$ws.Rows.Item(518).Insert()

# New row 518: Primera, 45147
$ws.Cells.Item(518, 1).Value = 7
$ws.Cells.Item(518, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(518, 3).Value = "Ñuble"
$ws.Cells.Item(518, 4).Value = 45147
$ws.Cells.Item(518, 5).Value = 16
$ws.Cells.Item(518, 6).Value = 100112009
$ws.Cells.Item(518, 7).Value = "Acelga"
$ws.Cells.Item(518, 8).Value = "Sin especificar"
$ws.Cells.Item(518, 9).Value = "Primera"
$ws.Cells.Item(518, 10).Value = 180
$ws.Cells.Item(518, 11).Value = 700
$ws.Cells.Item(518, 12).Value = 700
$ws.Cells.Item(518, 13).Value = 700
$ws.Cells.Item(518, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(518, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(518, 16).Value = 700
$ws.Cells.Item(518, 17).Value = 1
$ws.Cells.Item(518, 18).Value = "Hortaliza"

# New row 519: Segunda, 45147
$ws.Cells.Item(519, 1).Value = 7
$ws.Cells.Item(519, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(519, 3).Value = "Ñuble"
$ws.Cells.Item(519, 4).Value = 45147
$ws.Cells.Item(519, 5).Value = 16
$ws.Cells.Item(519, 6).Value = 100112009
$ws.Cells.Item(519, 7).Value = "Acelga"
$ws.Cells.Item(519, 8).Value = "Sin especificar"
$ws.Cells.Item(519, 9).Value = "Segunda"
$ws.Cells.Item(519, 10).Value = 150
$ws.Cells.Item(519, 11).Value = 500
$ws.Cells.Item(519, 12).Value = 500
$ws.Cells.Item(519, 13).Value = 500
$ws.Cells.Item(519, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(519, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(519, 16).Value = 500
$ws.Cells.Item(519, 17).Value = 1
$ws.Cells.Item(519, 18).Value = "Hortaliza"
